$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Parameter Values"

$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Results"
$ws2.Move($ws1, $null)

Write-Host $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
